$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C13 (row 23) is being replaced with a 4.7uF cap, same part used for C11/C2,C7,C9,C10,C19,C1,C3 (C1779)
$ws.Range("C23").Value = "4.7uF"
$ws.Range("D23").Value = "C1779"

# R14 (row 34) JLPCB part number reverts from shared-string "C31850" back to the same literal text
$ws.Range("D34").Value = "C31850"

# Move the active selection to D23, matching the edited cell
$ws.Range("D23").Select()
